$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.739.27'
$ws.Range('E2').Value = '  -2.60%  '

$ws.Range('D3').Value = '1.885.26'
$ws.Range('E3').Value = '  -5.13%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.57%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.06'
$ws.Range('E5').Value = '  -1.68%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.22%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4562'
$ws.Range('E7').Value = '  -1.78%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3789'
$ws.Range('E8').Value = '  -4.25%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.53'
$ws.Range('E9').Value = '  -1.55%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07707'
$ws.Range('E10').Value = '  -2.76%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9595'
$ws.Range('E11').Value = '  -4.31%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.92'
$ws.Range('E12').Value = '  -2.68%  '

$ws.Range('D13').Value = '1.885.10'
$ws.Range('E13').Value = '  -4.66%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.935'
$ws.Range('E14').Value = '  -3.97%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.639'
$ws.Range('E15').Value = '  -3.80%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06982'
$ws.Range('E16').Value = '  -1.78%  '

$ws.Range('E17').Value = '  +0.28%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '82.89'
$ws.Range('E18').Value = '  -6.71%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009474'
$ws.Range('E19').Value = '  -5.10%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.56'
$ws.Range('E20').Value = '  -3.82%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.34%  '

$ws.Range('D22').Value = '28.712.65'
$ws.Range('E22').Value = '  -2.98%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.301'
$ws.Range('E23').Value = '  -5.05%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.84'
$ws.Range('E24').Value = '  -3.76%  '

$ws.Range('D25').Value = '2.110.95'
$ws.Range('E25').Value = '  -4.83%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.077'
$ws.Range('E26').Value = '  -2.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.75'
$ws.Range('E27').Value = '  -2.00%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.94'
$ws.Range('E28').Value = '  -3.76%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.592'
$ws.Range('E29').Value = '  -7.07%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.78'
$ws.Range('E30').Value = '  -3.24%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.797'
$ws.Range('E31').Value = '  -6.51%  '

$ws.Range('E32').Value = '  -2.09%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8426'
$ws.Range('E33').Value = '  -5.86%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.049'
$ws.Range('E34').Value = '  -4.48%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.234'
$ws.Range('E35').Value = '  -8.55%  '

$ws.Range('E36').Value = '  -6.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05646'
$ws.Range('E37').Value = '  -3.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.139'
$ws.Range('E38').Value = '  -3.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.001'
$ws.Range('E39').Value = '  +0.40%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02023'
$ws.Range('E40').Value = '  -5.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.414'
$ws.Range('E41').Value = '  -6.63%  '

$ws.Range('E42').Value = '  -5.41%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.000002998'
$ws.Range('E43').Value = '  -25.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1742'
$ws.Range('E44').Value = '  -4.34%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.153'
$ws.Range('E45').Value = '  -6.80%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.689'
$ws.Range('E46').Value = '  +1.75%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5143'
$ws.Range('E47').Value = '  -4.33%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.26'
$ws.Range('E48').Value = '  -6.80%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06792'
$ws.Range('E49').Value = '  -2.83%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.056'
$ws.Range('E50').Value = '  -5.99%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.35'
$ws.Range('E51').Value = '  -2.51%  '
